$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 13.082
$ws.Range("E3").Value = 13.123
$ws.Range("E5").Value = 13.165
$ws.Range("A9").Value = -20.912
$ws.Range("E11").Value = 13.012
$ws.Range("E12").Value = 13
$ws.Range("A13").Value = -22.005
$ws.Range("A16").Value = -20.764
$ws.Range("A18").Value = -21.868
$ws.Range("A20").Value = -21.883
$ws.Range("E21").Value = 13.273
